# Rewrites the four "fake field" paragraphs (fldChar begin/instrText/fldChar end)
# into plain-text runs using the {m: ...} textual syntax, per
# TokenIteratorFieldRewriterSplit. We replace each paragraph's Range content
# wholesale via InsertXML with the equivalent WordprocessingML, preserving
# the surrounding run/bookmark structure exactly.

$d = $word.ActiveDocument

$pkgHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
$pkgFooter = '</w:document></pkg:xmlData></pkg:part></pkg:package>'

# --- Paragraph: {m: 2.myTemplate()}  ------------------------------------
$body2 = '<w:body><w:p>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">{m: 2.myTemplate()}</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">    </w:t></w:r>' +
  '<w:r><w:rPr><w:color w:val="FF0000"/><w:sz w:val="32"/><w:highlight w:val="lightGray"/></w:rPr><w:t>&lt;---</w:t></w:r>' +
  '<w:r><w:rPr><w:color w:val="FF0000"/><w:sz w:val="32"/><w:highlight w:val="lightGray"/></w:rPr><w:t>Couldn' + [char]0x27 + 't find the ' + [char]0x27 + 'myTemplate(java.lang.Integer)' + [char]0x27 + ' service</w:t></w:r>' +
  '</w:p></w:body>'

$p = $d.Paragraphs(2)
$p.Range.InsertXML($pkgHeader + $body2 + $pkgFooter)

# --- Paragraph: {m:template myTemplate(a:)}  (contains bookmark) --------
$body4 = '<w:body><w:p>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>{m:template myTemplate(a:</w:t></w:r>' +
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">)}</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">    </w:t></w:r>' +
  '<w:r><w:rPr><w:color w:val="FF0000"/><w:sz w:val="32"/><w:highlight w:val="lightGray"/></w:rPr><w:t>&lt;---</w:t></w:r>' +
  '<w:r><w:rPr><w:color w:val="FF0000"/><w:sz w:val="32"/><w:highlight w:val="lightGray"/></w:rPr><w:t>Expression &quot;a&quot; is invalid: missing type literal</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">    </w:t></w:r>' +
  '<w:r><w:rPr><w:color w:val="FF0000"/><w:sz w:val="32"/><w:highlight w:val="lightGray"/></w:rPr><w:t>&lt;---</w:t></w:r>' +
  '<w:r><w:rPr><w:color w:val="FF0000"/><w:sz w:val="32"/><w:highlight w:val="lightGray"/></w:rPr><w:t>missing type literal</w:t></w:r>' +
  '</w:p></w:body>'

$p = $d.Paragraphs(4)
$p.Range.InsertXML($pkgHeader + $body4 + $pkgFooter)

# --- Paragraph: {m: a + a}  ----------------------------------------------
$body5 = '<w:body><w:p>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">{m: a + a}</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">    </w:t></w:r>' +
  '<w:r><w:rPr><w:color w:val="FF0000"/><w:sz w:val="32"/><w:highlight w:val="lightGray"/></w:rPr><w:t>&lt;---</w:t></w:r>' +
  '<w:r><w:rPr><w:color w:val="FF0000"/><w:sz w:val="32"/><w:highlight w:val="lightGray"/></w:rPr><w:t>missing type literal</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">    </w:t></w:r>' +
  '<w:r><w:rPr><w:color w:val="FF0000"/><w:sz w:val="32"/><w:highlight w:val="lightGray"/></w:rPr><w:t>&lt;---</w:t></w:r>' +
  '<w:r><w:rPr><w:color w:val="FF0000"/><w:sz w:val="32"/><w:highlight w:val="lightGray"/></w:rPr><w:t>missing type literal</w:t></w:r>' +
  '</w:p></w:body>'

$p = $d.Paragraphs(5)
$p.Range.InsertXML($pkgHeader + $body5 + $pkgFooter)

# --- Paragraph: {m:endtemplate}  -----------------------------------------
$body6 = '<w:body><w:p>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">{m:endtemplate}</w:t></w:r>' +
  '</w:p></w:body>'

$p = $d.Paragraphs(6)
$p.Range.InsertXML($pkgHeader + $body6 + $pkgFooter)

Write-Host "edit applied"
